$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "BCT Cases For Item Module"
$ws.Range("A6").Value = "ItemModuleBCT"
$ws.Range("C6").Value = "Y"
